$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.83"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.76%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "19"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.69%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "19"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.085"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.07%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "19"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05704"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1.80%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "19"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.509"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.43%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "19"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8197"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.78%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "19"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8632"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2.03%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "19"
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1331"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.49%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "19"
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06936"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.92%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "19"
$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.02827"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.90%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "19"
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "BitMartToken"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09400"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.11%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "19"
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "BitForexToken"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.001533"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.24%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "19"
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = "CoinExToken"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.04071"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-12.54%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "19"
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "One"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005988"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-93.96%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "19"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006111"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.06%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "19"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.505"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-2.60%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "19"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.15%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "19"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.316"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "12.68%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "19"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "19"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.03167"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.12%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "19"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-1.76%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "19"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.558"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.91%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "19"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1373"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.71%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "19"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-2.55%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "19"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.003973"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-13.30%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "19"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.00009898"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.08%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "19"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001448"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "3.60%"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "19"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "19"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "19"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "19"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "19"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "19"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "19"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "19"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "19"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "19"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "19"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "19"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03723"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.45%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "19"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.005728"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-7.77%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "19"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.02%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "19"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002299"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-8.05%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "19"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009373"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "6.58%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "19"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005138"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-3.01%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "19"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.06%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "19"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-7.79%"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "19"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002529"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-4.51%"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "19"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.06%"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "19"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.06%"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "19"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "19"
